# COREESG_holdings.xlsx update
# - refresh the "Model holdings provided as of ..." date in the disclosure text (A10)
# - refresh the Weight (D) / Percent Change (E) figures for rows 2-7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect so the cells can be written, then restore
# protection afterwards.
$ws.Unprotect()

# Disclosure text: bump the "as of" date from 2021-03-30 to 2021-03-31
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-31 for illustrative purposes only and are subject to change."

# NULG (row 2)
$ws.Range("D2").Value = 0.2452741571500534
$ws.Range("E2").Value = 0.01665789935121853

# NULV (row 3)
$ws.Range("D3").Value = 0.4995218261319623
$ws.Range("E3").Value = -0.003303964757709221

# NUMG (row 4)
$ws.Range("D4").Value = 0.09798435107353365
$ws.Range("E4").Value = 0.0168562144597888

# NUMV (row 5)
$ws.Range("D5").Value = 0.09982288766282342
$ws.Range("E5").Value = -0.0005807200929152101

# NUSC (row 6)
$ws.Range("D6").Value = 0.05739677798162732
$ws.Range("E6").Value = 0.005608787099789625

# Total (row 7) - only the Percent Change column changes
$ws.Range("E7").Value = 0.004350952100727223

# Restore sheet protection to match the original document state.
$ws.Protect()
